$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($CellRef, $NewValue)
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = $origStyle
}

Set-TextCell "D2" "62.338.54"
Set-TextCell "E2" "  -5.48%  "
Set-TextCell "D3" "3.171.47"
Set-TextCell "E3" "  -6.60%  "
Set-TextCell "E4" "  -0.03%  "
Set-TextCell "B5" "Solana"
Set-TextCell "C5" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D5" "171.30"
Set-TextCell "E5" "  -9.06%  "
Set-TextCell "B6" "BNB"
Set-TextCell "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell "D6" "507.83"
Set-TextCell "E6" "  -4.15%  "
Set-TextCell "D7" "0.583"
Set-TextCell "E8" "  +0.14%  "
Set-TextCell "D9" "3.168.63"
Set-TextCell "E9" "  -6.63%  "
Set-TextCell "D10" "0.588"
Set-TextCell "E10" "  -6.87%  "
Set-TextCell "D11" "51.63"
Set-TextCell "E11" "  -11.28%  "
Set-TextCell "D12" "0.126"
Set-TextCell "E12" "  -7.68%  "
Set-TextCell "E13" "  -4.03%  "
Set-TextCell "D14" "8.73"
Set-TextCell "E14" "  -7.23%  "
Set-TextCell "D15" "3.699.53"
Set-TextCell "E15" "  -6.15%  "
Set-TextCell "D16" "3.178.88"
Set-TextCell "E16" "  -6.28%  "
Set-TextCell "E17" "  -8.13%  "
Set-TextCell "D18" "62.298.19"
Set-TextCell "E18" "  -5.25%  "
Set-TextCell "D19" "16.84"
Set-TextCell "E19" "  -4.29%  "
Set-TextCell "D20" "10.69"
Set-TextCell "E20" "  -5.72%  "
Set-TextCell "D21" "0.936"
Set-TextCell "E21" "  -5.17%  "
Set-TextCell "D22" "356.93"
Set-TextCell "E22" "  -7.22%  "
Set-TextCell "B23" "PancakeSwap"
Set-TextCell "C23" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D23" "3.63"
Set-TextCell "E23" "  -3.43%  "
Set-TextCell "B24" "Litecoin"
Set-TextCell "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D24" "78.95"
Set-TextCell "E24" "  -5.59%  "
Set-TextCell "B25" "RenderToken"
Set-TextCell "C25" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D25" "10.86"
Set-TextCell "E25" "  -3.16%  "
Set-TextCell "D26" "5.93"
Set-TextCell "E26" "  -1.87%  "
Set-TextCell "E27" "  +0.13%  "
Set-TextCell "D28" "2.55"
Set-TextCell "E28" "  -5.28%  "
Set-TextCell "D29" "10.96"
Set-TextCell "E29" "  -5.48%  "
Set-TextCell "D30" "7.97"
Set-TextCell "E30" "  -6.63%  "
Set-TextCell "D31" "641.30"
Set-TextCell "E31" "  -6.34%  "
Set-TextCell "D32" "27.72"
Set-TextCell "E32" "  -7.46%  "
Set-TextCell "D33" "6.13"
Set-TextCell "E33" "  -9.01%  "
Set-TextCell "D34" "10.88"
Set-TextCell "E34" "  -3.55%  "
Set-TextCell "E35" "  -4.78%  "
Set-TextCell "D36" "56.76"
Set-TextCell "E36" "  -8.30%  "
Set-TextCell "E37" "  -0.10%  "
Set-TextCell "D38" "35.69"
Set-TextCell "E38" "  -3.35%  "
Set-TextCell "E39" "  -3.48%  "
Set-TextCell "E40" "  +0.07%  "
Set-TextCell "D41" "0.0₃0673"
Set-TextCell "E41" "  +6.55%  "
Set-TextCell "B42" "Kaspa"
Set-TextCell "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D42" "0.120"
Set-TextCell "E42" "  -4.55%  "
Set-TextCell "B43" "Maker"
Set-TextCell "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D43" "2.808.26"
Set-TextCell "E43" "  -4.00%  "
Set-TextCell "D44" "2.47"
Set-TextCell "E44" "  +2.63%  "
Set-TextCell "D45" "2.67"
Set-TextCell "E45" "  +0.93%  "
Set-TextCell "B46" "VeChain"
Set-TextCell "C46" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D46" "0.0379"
Set-TextCell "E46" "  -2.56%  "
Set-TextCell "B47" "ThetaToken"
Set-TextCell "C47" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D47" "2.52"
Set-TextCell "E47" "  -10.37%  "
Set-TextCell "B48" "Stacks"
Set-TextCell "C48" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D48" "2.70"
Set-TextCell "E48" "  +5.06%  "
Set-TextCell "B49" "Monero"
Set-TextCell "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D49" "134.12"
Set-TextCell "E49" "  -0.43%  "
Set-TextCell "B50" "Stellar"
Set-TextCell "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D50" "0.121"
Set-TextCell "E50" "  -5.28%  "
Set-TextCell "B51" "ApeXProtocol"
Set-TextCell "C51" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell "D51" "2.82"
Set-TextCell "E51" "  -3.86%  "
